$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ativação: date-like text "01/01/2021" ---------------------------------
# Assigning a date-shaped string straight into B8/C8 would make Excel
# reinterpret it as a real date (new number format + style). Write it into a
# scratch cell first, force that scratch cell to Text format so the literal
# string sticks, then copy *values only* into the target cells so the
# destination keeps its original style/number format untouched. Finally
# remove the scratch column so the sheet's used range is unchanged.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "01/01/2021"
$ws.Range("Z1").Copy()
$ws.Range("B8").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("C8").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").EntireColumn.Delete()

# --- Objetivos: trimmed to first sentence, added a space after leading "-" -
$objetivos = "- Apresentar a evolução das condições geológicas da Terra que culminaram com os recursos naturais existentes hoje, com ênfase nas reservas de combustíveis fósseis, hídricos e de minérios e a conseqüente reserva de energia advinda dessas fontes naturais. A América do Sul e do Brasil, mais especificamente, são destacados."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# --- Programa resumido: rewritten, shorter ---------------------------------
$programaResumido = "- Desenvolvimento da Terra.- Recursos minerais.- Matérias-primas da grande indústria metalúrgica: metais ferrosos e não-ferrosos"
$ws.Range("B14").Value = $programaResumido
$ws.Range("C14").Value = $programaResumido

# --- Programa: rewritten, shorter ------------------------------------------
$programa = "- Desenvolvimento da Terra. - Principais Eras Geológicas. - Matérias-primas para a grande indústria metalúrgica: metais ferrosos e metais não-ferrosos."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# --- Bibliografia: rewritten reference list ---------------------------------
$bibliografia = "- MILLER, Jr. G. T. “Ciência Ambiental”,  Editora: Thomson, 2006.- ABREU, S. F. “Recursos Minerais do Brasil”, Editora: Edgard Blücher, 1973.-  SKINNER, B. J. “Recursos Minerais da Terra”, Editora: Edgard Blücher, 1996.- WICANDER, R.; MONROE, J. S. “Fundamentos de Geologia”, Editora: Cengage Learning, 2009. - PRESS, F.; Siever, R.; Jordan, T.; Grotzinger, J. “Para Entender a Terra”, Editora: Bookman,  2006.- SCHÄFER, A. “Fundamentos de Ecologia e Biogeografia de Águas Continentais”, Editora: Universidade, Porto Alegre. - Revistas especializadas e sites, dado ao caráter dinâmico das informações sobre reservas minerais e recursos naturais em geral."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia
